# Appends the new "Responsive design" section (14 bulleted paragraphs plus a
# trailing blank paragraph) to the end of the document body, right before </w:body>
# / <w:sectPr>, matching the target diff exactly.
#
# Strategy: Range.InsertXML() is the most reliable way to land paragraphs with
# exact run/formatting fidelity (bCs/szCs, mixed bold runs, <w:lastRenderedPageBreak/>,
# xml:space="preserve", literal "<"/">" via &lt;/&gt;, etc.) verbatim. InsertXML always
# consumes the *text* of the target range but leaves that range's own paragraph mark
# behind as a new, now-empty, trailing paragraph -- so we chain inserts by always
# targeting $d.Paragraphs.Last. Doing this at the literal end-of-story position can
# clobber the preceding paragraph, so we first seed a placeholder trailing paragraph
# via InsertParagraphAfter() and insert into that instead.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$frag0 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Responsive design</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag0) | Out-Null

$frag1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Mobile-first VS Desktop-first and breakpoints</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag1) | Out-Null

$frag2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">Desktop fist: Start writing CSS for the desktop: large screen; - Then, media queries shrink design to smaller screens. </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag2) | Out-Null

$frag3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Mobile first: Start writing CSS for mobile devices: small screen; Then, media queries expand design to a large desktop screen; -Forces us to reduce websites and apps to the absolute essentials.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag3) | Out-Null

$frag4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Media queries don’t add any importance or specificity to selectors, so code order matters – media queries at the end.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag4) | Out-Null

$frag5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">Is Mobile-first right for you? </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>PROS</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">: -100% optimized for the mobile experience; -Reduces websites and apps to the absolute essentials; -Results in smaller, faster and more efficient products; -Prioritizes content over aesthetic design, which may be desirable. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>CONS</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>: -The desktop version might feel overly empty and simplistic; -More difficult and counterintuitive to develop; -Less creative freedom, makin ti more difficult to create distinctive products; -Clients are used to see a desktop version of the site as a prototype; -Do your users even use the mobile internet? What’s the purpose of your website?</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag5) | Out-Null

$frag6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>No matter what you do, always keeps both desktop and mobile in mind</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag6) | Out-Null

$frag7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Selecting the breakpoints: The options: Bad(devices screens-like iphone screens), Good(group different devices-mobiles, tablets, laptops), Perfect(content into your design, design for some screen and when the content breaks define a new break point).</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag7) | Out-Null

$frag8 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>How to use a powerful Sass mixing to write all our media queries; -How to use the @content and @if Sass directives; -Taking advantage of Chrome DevTools for responsive design.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag8) | Out-Null

$frag9 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Responsive Images</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>: The goal of responsive images is to serve the right image to the right screen size and device, in order to avoid downloading unnecessary large images on smaller screens.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag9) | Out-Null

$frag10 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>3 Uses cases: Resolution switching(decrease image resolution on smaller screen), Density switching-special case of resolution switching but the screen size doesn’t matter(double the resolution on the photo), Art direction(other image, remodeled)</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag10) | Out-Null

$frag11 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Responsive images in HTML – how to use the srcset attribute on the &lt;img&gt; and &lt;source&gt; elements, together with density descriptions; -how and why to use the &lt;picture&gt; element for art direction; -how to write media queries in HTML.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>-How to allow the browser to decide the best image to download, using the srcset attribute, width descriptors, and the sizes attribute of the &lt;img&gt; element.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag11) | Out-Null

$frag12 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Responsive images in CSS: -how to implement responsive images in CSS; -How to use resolution media queries to target high-resolution screens with 2x; How to combine multiple conditions in media queries.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag12) | Out-Null

$frag13 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="12"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Browser support: Always check caniuse.com before using a Modern css property in production; -use graceful degradation with @supports</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>. -How to use @supports feature queries; -Implement graceful degradation on selected properties; -How to use backdrop-filter</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>(sepia, brightness, blur…</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> applied on the background elements</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($frag13) | Out-Null

# Final trailing empty paragraph (ind left=360, no list numbering / pStyle).
$fragFinal = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target = $d.Paragraphs.Last
$target.Range.InsertXML($fragFinal) | Out-Null

# InsertXML left the consumed paragraph's own end-of-paragraph mark behind as one
# more (now spurious) empty ListParagraph-styled paragraph; merge it away so exactly
# one final paragraph (the ind left=360 one) remains before <w:sectPr>.
$prevPara = $d.Paragraphs[$d.Paragraphs.Count - 1]
$extraPara = $d.Paragraphs.Last
$mergeRange = $d.Range($prevPara.Range.End - 1, $extraPara.Range.End)
$mergeRange.Delete() | Out-Null

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
